$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.857.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.446.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.05%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -2.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.443.24"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.69%  "
$ws.Range("E10").Value = "  -6.53%  "
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("E12").Value = "  -6.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.884.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.95%  "
$ws.Range("E16").Value = "  -7.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.702.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.446.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.01%  "
$ws.Range("E21").Value = "  -6.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "317.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.78%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.33%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0₃0966"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.04%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.573.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "537.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("E33").Value = "  -6.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.66%  "
$ws.Range("E35").Value = "  -7.66%  "
$ws.Range("E36").Value = "  -9.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.376"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.30%  "
$ws.Range("E41").Value = "  -6.32%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.02%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.01%  "
$ws.Range("E45").Value = "  -6.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "145.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.80%  "
$ws.Range("E47").Value = "  -6.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0527"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.20%  "
$ws.Range("E50").Value = "  -6.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0937"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.20%  "
